$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

# Update r_s_star (J2) and h_p_star (K2) values for the 15mmol_20C test
$ws.Range("J2").Value = 0.0104
$ws.Range("K2").Value = 0.2732
